# Refactoring 9/26/24 @ 14:53
#
# Adds three new worksheets (RequestLogin, Data, Session) after the existing
# RequestSignup sheet and populates each with a header row + one sample data
# row, mirroring the login / user-data / session payloads used by the app.

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd($name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $last)
    $newSheet.Name = $name
    return $newSheet
}

# ---------------------------------------------------------------------------
# 1. New sheets, in order: RequestLogin, Data, Session
# ---------------------------------------------------------------------------
$wsLogin   = Add-SheetAtEnd "RequestLogin"
$wsData    = Add-SheetAtEnd "Data"
$wsSession = Add-SheetAtEnd "Session"

# ---------------------------------------------------------------------------
# 2. RequestLogin sheet
# ---------------------------------------------------------------------------
$wsLogin.Cells.Item(1, 1).Value = "username"
$wsLogin.Cells.Item(1, 2).Value = "password"
$wsLogin.Cells.Item(1, 3).Value = "longitude"
$wsLogin.Cells.Item(1, 4).Value = "latitude"
$wsLogin.Cells.Item(1, 5).Value = "deviceToken"
$wsLogin.Cells.Item(1, 6).Value = "userAgent"
$wsLogin.Cells.Item(1, 7).Value = "ip"

$wsLogin.Cells.Item(2, 1).Value = "testuser40"
$wsLogin.Cells.Item(2, 2).Value = 'Test1234$'
$wsLogin.Cells.Item(2, 3).Value = 1234567.0
$wsLogin.Cells.Item(2, 4).Value = 1234567.0
$wsLogin.Cells.Item(2, 5).Value = "deviceToken"
$wsLogin.Cells.Item(2, 6).Value = "useragent"
$wsLogin.Cells.Item(2, 7).Value = "ip"

# ---------------------------------------------------------------------------
# 3. Data sheet
# ---------------------------------------------------------------------------
$wsData.Cells.Item(1, 1).Value  = "id"
$wsData.Cells.Item(1, 2).Value  = "email"
$wsData.Cells.Item(1, 3).Value  = "username"
$wsData.Cells.Item(1, 4).Value  = "password"
$wsData.Cells.Item(1, 5).Value  = "firstName"
$wsData.Cells.Item(1, 6).Value  = "lastName"
$wsData.Cells.Item(1, 7).Value  = "phone"
$wsData.Cells.Item(1, 8).Value  = "providerGivenID"
$wsData.Cells.Item(1, 9).Value  = "confirmationCode"
$wsData.Cells.Item(1, 10).Value = "roles"
$wsData.Cells.Item(1, 11).Value = "isAuthenticated"
$wsData.Cells.Item(1, 12).Value = "isAccountNonExpired"
$wsData.Cells.Item(1, 13).Value = "isAccountNonLocked"
$wsData.Cells.Item(1, 14).Value = "isCredentialsNonExpired"
$wsData.Cells.Item(1, 15).Value = "isEnabled"

$wsData.Cells.Item(2, 1).Value  = 199.0
$wsData.Cells.Item(2, 2).Value  = "testuser40@gmail.com"
$wsData.Cells.Item(2, 3).Value  = "testuser40"
$wsData.Cells.Item(2, 4).Value  = 'Test1234$'
$wsData.Cells.Item(2, 5).Value  = "Test"
$wsData.Cells.Item(2, 6).Value  = "One"
$wsData.Cells.Item(2, 7).Formula = "=14048205065"
$wsData.Cells.Item(2, 8).Value  = "jeWAlxrTlRRRxkQbpyPBgt2Ogl72"
$wsData.Cells.Item(2, 9).Value  = 111111111111
$wsData.Cells.Item(2, 10).Value = 'a:1:{s:13:"administrator";b:1;}'''
$wsData.Cells.Item(2, 11).Value = 1.0
$wsData.Cells.Item(2, 12).Value = 1.0
$wsData.Cells.Item(2, 13).Value = 1.0
$wsData.Cells.Item(2, 14).Value = 1.0
$wsData.Cells.Item(2, 15).Value = 1.0

# ---------------------------------------------------------------------------
# 4. Session sheet (Data columns A-O + login/session columns P-V)
# ---------------------------------------------------------------------------
$wsSession.Cells.Item(1, 1).Value  = "id"
$wsSession.Cells.Item(1, 2).Value  = "email"
$wsSession.Cells.Item(1, 3).Value  = "username"
$wsSession.Cells.Item(1, 4).Value  = "password"
$wsSession.Cells.Item(1, 5).Value  = "firstName"
$wsSession.Cells.Item(1, 6).Value  = "lastName"
$wsSession.Cells.Item(1, 7).Value  = "phone"
$wsSession.Cells.Item(1, 8).Value  = "providerGivenID"
$wsSession.Cells.Item(1, 9).Value  = "confirmationCode"
$wsSession.Cells.Item(1, 10).Value = "roles"
$wsSession.Cells.Item(1, 11).Value = "isAuthenticated"
$wsSession.Cells.Item(1, 12).Value = "isAccountNonExpired"
$wsSession.Cells.Item(1, 13).Value = "isAccountNonLocked"
$wsSession.Cells.Item(1, 14).Value = "isCredentialsNonExpired"
$wsSession.Cells.Item(1, 15).Value = "isEnabled"
$wsSession.Cells.Item(1, 16).Value = "longitude"
$wsSession.Cells.Item(1, 17).Value = "latitude"
$wsSession.Cells.Item(1, 18).Value = "deviceToken"
$wsSession.Cells.Item(1, 19).Value = "userAgent"
$wsSession.Cells.Item(1, 20).Value = "ip"
$wsSession.Cells.Item(1, 21).Value = "accessToken"
$wsSession.Cells.Item(1, 22).Value = "refreshToken"

$wsSession.Cells.Item(2, 1).Value  = 199.0
$wsSession.Cells.Item(2, 2).Value  = "testuser40@gmail.com"
$wsSession.Cells.Item(2, 3).Value  = "testuser40"
$wsSession.Cells.Item(2, 4).Value  = 'Test1234$'
$wsSession.Cells.Item(2, 5).Value  = "Test"
$wsSession.Cells.Item(2, 6).Value  = "One"
$wsSession.Cells.Item(2, 7).Formula = "=14048205065"
$wsSession.Cells.Item(2, 8).Value  = "jeWAlxrTlRRRxkQbpyPBgt2Ogl72"
$wsSession.Cells.Item(2, 9).Value  = 111111111111
$wsSession.Cells.Item(2, 10).Value = 'a:1:{s:13:"administrator";b:1;}'''
$wsSession.Cells.Item(2, 11).Value = 1.0
$wsSession.Cells.Item(2, 12).Value = 1.0
$wsSession.Cells.Item(2, 13).Value = 1.0
$wsSession.Cells.Item(2, 14).Value = 1.0
$wsSession.Cells.Item(2, 15).Value = 1.0
$wsSession.Cells.Item(2, 16).Value = 123456.0
$wsSession.Cells.Item(2, 17).Value = 123456.0
$wsSession.Cells.Item(2, 18).Value = "deviceToken"
$wsSession.Cells.Item(2, 19).Value = "userAgent"
$wsSession.Cells.Item(2, 20).Value = "123.456.7891"
$wsSession.Cells.Item(2, 21).Value = "eyJhbGciOiJIUzI1NiJ9.eyJsb2NhdGlvbiI6eyJsb25naXR1ZGUiOjEyMzQ1NjcuMCwibGF0aXR1ZGUiOjEyMzQ1NjcuMH0sImlzcyI6Im9yYi1nYXRld2F5Iiwic3ViIjoidGVzdHVzZXI0MCIsImlhdCI6MTcyNzI5NzA5MCwiZXhwIjoxNzI3MzgzNDkwfQ._fdNIP_Bl-S3MNaVlxg4yFSi3yOoxaf7IaCksGqY4ak"
$wsSession.Cells.Item(2, 22).Value = "eyJhbGciOiJIUzI1NiJ9.eyJpc3MiOiJvcmItZ2F0ZXdheSIsInN1YiI6InRlc3R1c2VyNDAiLCJpYXQiOjE3MjcyOTcwOTAsImV4cCI6MTcyNzM4MzQ5MH0.5oiiLwWopT5SOBRLy208Oi4gceYeOhnOphs7hm6xU8Q"

# Leave RequestSignup (sheet1) as the active sheet, matching the original file.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
